$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the printed date in A1 by one day (2024-01-17 -> 2024-01-18,
# serial 45308 -> 45309)
$ws.Range("A1").Value2 = 45309

# Step 2: update the price in D29 (960 -> 1570)
$ws.Range("D29").Value2 = 1570
